$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('TOLENTINO VASQUEZ DIANA KATHERYN', 60),
    @('JULCA VALENZUELA CINTIA KARYN', 57),
    @('CARRILLO MARTÍNEZ HEIDY NAYELI', 51),
    @('DE LA CRUZ BENITES RICHARD ALEXANDER', 50),
    @('YZQUIERDO CARHUATANTA LEYDY YANELA', 48),
    @('RODRIGUEZ RUBIO SANDRA MABEL', 48),
    @('ARENAS ZAVALA ANDYELA PATRICIA ISIDORA', 45),
    @('REYES RODRIGUEZ JEISSON STEVEN', 44),
    @('GASLAC GUTIERREZ FRANK JHORDY', 43),
    @('SANCHEZ CORTEZ LEYLA DIANA', 43),
    @('RUBIO MARIÑOS GISELA JUDITH', 42),
    @('VALER VEGA PATRICIA GERALDINE', 41),
    @('PIERINA NAGIELLY SANDOVAL CONTRERAS', 41),
    @('CYNTHIA RODRIGUEZ LECCA', 40),
    @('PONCE VILLANUEVA CARMEN ISABEL', 36),
    @('GUZMAN ZAVALETA CECILIA MARISOL', 33),
    @('SEGURA ASTO YAMILET ANTONELA', 25),
    @('RODRIGUEZ VASQUEZ WALTER', 20),
    @('LEON VERA MELISSA FIORELLA', 16),
    @('RUTH MELISSA RAMIREZ VELEZMORO', 16),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

Write-Output "Updated $($data.Count) rows"
